$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (price + 1h volume change) scraped for this run.
# Column D (Price) values are forced to Text so values such as "12.30" or
# "0.0740" keep their exact original formatting instead of being silently
# normalised to a number (12.3 / 0.074) by the Value setter.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.629.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.260.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.70%  "
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.55"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.611.29"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.273.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.757"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.561.48"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +10.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0899"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  +4.47%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +5.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.99"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.10"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.89%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.12"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0740"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.24%  "
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("E38").Value = "  +2.65%  "
$ws.Range("E39").Value = "  +7.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.057.27"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.24"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0276"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("E47").Value = "  +5.32%  "
$ws.Range("E48").Value = "  +3.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.37%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.75%  "
